# status report.docx edits
# 1. Merge "05/16" + "/2019" into a single run "05/16/2019" (table 1, date cell)
# 2. Henry row: bug fixes 6 -> 9 hours; Actual Hours 19 -> 22
# 3. Harman row: bug fixes 8 -> 9 hours; Actual Hours 23 -> 24
# 4. Oliver row: bug fixes 8 -> 10 hours; Actual Hours 23 -> 25
# 5. John row: bug fixes 4 -> 5 hours; Actual Hours 27 -> 28
# 6. Merge "QR code generation, tags for AI" + "  " into a single run (no visible text change)

$d = $word.ActiveDocument

# --- 1. Date cell: merge "05/16/2019" across the two runs into one run ---
# (search text is unique document-wide, so a document-scoped Find is safe)
$d.Content.Find.Execute(
    "05/16/2019", $true, $false, $false, $false, $false, $true, 1, $false,
    "05/16/2019", 2) | Out-Null

# --- Table 2: "Tasks Completed This Period" ---
$t2 = $d.Tables.Item(2)

# Henry (row 2): task description bug-fixes hours 6 -> 9 (unique phrase document-wide)
$d.Content.Find.Execute(
    "Testing User, Admin, and Business features (13 hours), bug fixes (6 hours)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Testing User, Admin, and Business features (13 hours), bug fixes (9 hours)", 2) | Out-Null
# Henry Actual Hours 19 -> 22 (set the cell range text directly; plain numbers are not
# unique document-wide, so Find must not be used here)
$t2.Cell(2, 4).Range.Text = "22"

# Harman (row 3): task description bug-fixes hours 8 -> 9 (unique phrase document-wide)
$d.Content.Find.Execute(
    "Testing User, Admin, and Business features (11 hours), bug fixes (8 hours), finish user guide (4 hrs)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Testing User, Admin, and Business features (11 hours), bug fixes (9 hours), finish user guide (4 hrs)", 2) | Out-Null
# Harman Actual Hours 23 -> 24
$t2.Cell(3, 4).Range.Text = "24"

# Oliver (row 4): task description bug-fixes hours 8 -> 10 (unique phrase document-wide)
$d.Content.Find.Execute(
    "Testing User, Admin, and Business features (10 hours), bug fixes (8 hours), user documentation formatting and revision (5 hrs)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Testing User, Admin, and Business features (10 hours), bug fixes (10 hours), user documentation formatting and revision (5 hrs)", 2) | Out-Null
# Oliver Actual Hours 23 -> 25
$t2.Cell(4, 4).Range.Text = "25"

# John (row 5): task description bug-fixes hours 4 -> 5 (unique phrase document-wide)
$d.Content.Find.Execute(
    "Testing User, Admin, and Business features (6 hours), bug fixes (4 hours), UI unification (12 hours), final presentation (5 hours)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Testing User, Admin, and Business features (6 hours), bug fixes (5 hours), UI unification (12 hours), final presentation (5 hours)", 2) | Out-Null
# John Actual Hours 27 -> 28
$t2.Cell(5, 4).Range.Text = "28"

# --- Table 3: "Tasks Planned but Not Completed" ---
$t3 = $d.Tables.Item(3)

# Harman (row 2): merge "QR code generation, tags for AI" + trailing two spaces into one run
$t3.Cell(2, 2).Range.Text = "QR code generation, tags for AI  "
